$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.47"
$ws.Range("E2").Value = "'5.62%"
$ws.Range("D3").Value = "'32.00"
$ws.Range("E3").Value = "'9.88%"
$ws.Range("D4").Value = "'5.268"
$ws.Range("E4").Value = "'1.24%"
$ws.Range("D5").Value = "'0.07479"
$ws.Range("E5").Value = "'7.43%"
$ws.Range("D6").Value = "'7.836"
$ws.Range("D7").Value = "'3.809"
$ws.Range("E7").Value = "'7.15%"
$ws.Range("E8").Value = "'9.05%"
$ws.Range("D9").Value = "'0.9180"
$ws.Range("E9").Value = "'2.09%"
$ws.Range("D10").Value = "'0.1679"
$ws.Range("E10").Value = "'4.62%"
$ws.Range("D11").Value = "'0.08019"
$ws.Range("E11").Value = "'6.55%"
$ws.Range("D12").Value = "'0.08034"
$ws.Range("E12").Value = "'3.90%"
$ws.Range("D13").Value = "'0.03004"
$ws.Range("E13").Value = "'2.30%"
$ws.Range("D14").Value = "'0.09885"
$ws.Range("E14").Value = "'9.72%"
$ws.Range("D15").Value = "'0.001491"
$ws.Range("E15").Value = "'-5.14%"
$ws.Range("D16").Value = "'0.04598"
$ws.Range("E16").Value = "'1.50%"
$ws.Range("D17").Value = "'0.006567"
$ws.Range("E17").Value = "'0.46%"
$ws.Range("D18").Value = "'3.468"
$ws.Range("E18").Value = "'-0.52%"
$ws.Range("D20").Value = "'0.3326"
$ws.Range("E20").Value = "'2.50%"
$ws.Range("E21").Value = "'0.10%"
$ws.Range("D22").Value = "'4.485"
$ws.Range("E22").Value = "'10.56%"
$ws.Range("D23").Value = "'0.1620"
$ws.Range("E23").Value = "'1.27%"
$ws.Range("D24").Value = "'0.001216"
$ws.Range("E24").Value = "'0.55%"
$ws.Range("D25").Value = "'0.004447"
$ws.Range("E25").Value = "'7.38%"
$ws.Range("D26").Value = "'0.0001398"
$ws.Range("E26").Value = "'19.56%"
$ws.Range("D27").Value = "'0.0001775"
$ws.Range("E27").Value = "'6.09%"
$ws.Range("D39").Value = "'0.01724"
$ws.Range("E39").Value = "'2,548.84%"
$ws.Range("D40").Value = "'0.04498"
$ws.Range("E40").Value = "'3.07%"
$ws.Range("D41").Value = "'0.007153"
$ws.Range("E41").Value = "'3.27%"
$ws.Range("D42").Value = "'0.1349"
$ws.Range("E42").Value = "'8.20%"
$ws.Range("D43").Value = "'0.002157"
$ws.Range("E43").Value = "'4.26%"
$ws.Range("D44").Value = "'0.01283"
$ws.Range("E44").Value = "'8.79%"
$ws.Range("D45").Value = "'0.00006176"
$ws.Range("E45").Value = "'5.95%"
$ws.Range("D46").Value = "'1.868"
$ws.Range("E46").Value = "'-3.17%"
$ws.Range("D47").Value = "'0.01298"
$ws.Range("E47").Value = "'-0.13%"
